$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("L2").Value = 'https://vehicle-images.dealerinspire.com/59aa-110010848/1FBAX2C84NKA19843/4c9e05ac2920a65f2677ba6b82038fc0.png, https://vehicle-images.dealerinspire.com/fe93-110010848/1FBAX2C84NKA19843/072c767a5bd51f2e08b53b840ccb0e26.jpg, https://vehicle-images.dealerinspire.com/4c57-110010848/1FBAX2C84NKA19843/e0ac005380209e33d30655657f9e66a2.jpg, https://vehicle-images.dealerinspire.com/e8ca-110010848/1FBAX2C84NKA19843/bd6db14aae01632f809db9193fc7555e.jpg, https://vehicle-images.dealerinspire.com/2e6a-110010848/1FBAX2C84NKA19843/6c991a50868af59b825df0455938439f.jpg, https://vehicle-images.dealerinspire.com/cdb8-110010848/1FBAX2C84NKA19843/88dcce73e3a7bbbda8d9be8b2580ec44.jpg, https://vehicle-images.dealerinspire.com/41e7-110010848/1FBAX2C84NKA19843/4395f3a3cf9cf892acdabcef1cd5bd40.jpg, https://vehicle-images.dealerinspire.com/6d1f-110010848/1FBAX2C84NKA19843/370410c51134e42814ce09df2791691b.jpg, https://vehicle-images.dealerinspire.com/84b0-110010848/1FBAX2C84NKA19843/d80ed6de17ebf61e851fe563b7c43117.jpg, https://vehicle-images.dealerinspire.com/5084-110010848/1FBAX2C84NKA19843/12a81fd63918b8ada097ec17c453af60.jpg, https://vehicle-images.dealerinspire.com/a49a-110010848/1FBAX2C84NKA19843/bbc96781c3ec88c7f32cd49b97cc2864.jpg, https://vehicle-images.dealerinspire.com/b0be-110010848/1FBAX2C84NKA19843/4ca962a02414c336f06f49d725873c24.jpg, https://vehicle-images.dealerinspire.com/389f-110010848/1FBAX2C84NKA19843/30fff4b3d8b43c35163b4a619c31134c.jpg, https://vehicle-images.dealerinspire.com/6742-110010848/1FBAX2C84NKA19843/0dd19f0576fd6f6be19329bd19c71050.jpg, https://vehicle-images.dealerinspire.com/5304-110010848/1FBAX2C84NKA19843/f969abbccd1361aabe79e19cd62308b8.jpg, https://vehicle-images.dealerinspire.com/c31c-110010848/1FBAX2C84NKA19843/9d1f27fa05b25d05e9faffc68a9e95ef.jpg, https://vehicle-images.dealerinspire.com/0192-110010848/1FBAX2C84NKA19843/2a1b0096971d1649e65601780d3566fd.jpg, https://vehicle-images.dealerinspire.com/1a9b-110010848/1FBAX2C84NKA19843/e4c91d668662b20f55cacdad36c9b941.jpg, https://vehicle-images.dealerinspire.com/8bfa-110010848/1FBAX2C84NKA19843/79efffb0bfb8aae69bc681a52540322f.jpg, https://vehicle-images.dealerinspire.com/6b8f-110010848/1FBAX2C84NKA19843/c81216f8dcf004512cdd37e38025cf75.jpg, https://vehicle-images.dealerinspire.com/a9db-110010848/1FBAX2C84NKA19843/ca1faf5198f57dd697258305350f726a.jpg, https://vehicle-images.dealerinspire.com/dc77-110010848/1FBAX2C84NKA19843/03f03ab0f70df762c8f53e665fb0a798.jpg, https://vehicle-images.dealerinspire.com/7025-110010848/1FBAX2C84NKA19843/ce8f2b90ccb941e2bb7418c9cf6d165b.jpg, https://vehicle-images.dealerinspire.com/5d18-110010848/1FBAX2C84NKA19843/ed04c810a2a5f7ca59302ad20a90eb9c.jpg, https://vehicle-images.dealerinspire.com/c313-110010848/1FBAX2C84NKA19843/5044786422e7bccca0f5e9925143ee31.jpg, https://vehicle-images.dealerinspire.com/cc1d-110010848/1FBAX2C84NKA19843/44eaa74e686f8853afe262f58d90046e.jpg, https://vehicle-images.dealerinspire.com/04bb-110010848/1FBAX2C84NKA19843/882a183f5b741b37a201624cc96e70e8.jpg, https://vehicle-images.dealerinspire.com/ca35-110010848/1FBAX2C84NKA19843/beb56e55770d9d2dea0c58eb1a6ef708.jpg'
$ws.Range("L3").Value = 'https://vehicle-images.dealerinspire.com/5d88-11001792/W1Z70FGYXPT140876/9e61bebbc03cc097b80de804d2557d49.jpg, https://vehicle-images.dealerinspire.com/6b67-11001792/W1Z70FGYXPT140876/95eebcb5a959d8f587c65f7ef9b303d9.jpg, https://vehicle-images.dealerinspire.com/3033-11001792/W1Z70FGYXPT140876/aee0b823893010062b2abef3dc4aac94.jpg, https://vehicle-images.dealerinspire.com/230c-11001792/W1Z70FGYXPT140876/53fab97009aa47449f3d1699b3daf360.jpg, https://vehicle-images.dealerinspire.com/354f-11001792/W1Z70FGYXPT140876/d4bcda1d246989e8c05931f85d4b84ce.jpg, https://vehicle-images.dealerinspire.com/a3ae-11001792/W1Z70FGYXPT140876/b1d968620ef1ce173f7c04d6c0e1821a.jpg, https://vehicle-images.dealerinspire.com/f62e-11001792/W1Z70FGYXPT140876/75849829be54b344cd1c274352783c35.jpg, https://vehicle-images.dealerinspire.com/2a80-11001792/W1Z70FGYXPT140876/4aea4eecf2d918f159924abe184ce033.jpg, https://vehicle-images.dealerinspire.com/3cd4-11001792/W1Z70FGYXPT140876/a77216eddeac6a3acda0f8e54ec4d73a.jpg, https://vehicle-images.dealerinspire.com/96ac-11001792/W1Z70FGYXPT140876/a353ff22082bf6e746dd937606560a10.jpg, https://vehicle-images.dealerinspire.com/ba40-11001792/W1Z70FGYXPT140876/cb54c1932b3c431f23be704372211c18.jpg, https://vehicle-images.dealerinspire.com/d8a7-11001792/W1Z70FGYXPT140876/2120b223bcc2f682d513af15d0af660d.jpg, https://vehicle-images.dealerinspire.com/551f-11001792/W1Z70FGYXPT140876/b61142981f19782c47396c5e71748fbf.jpg, https://vehicle-images.dealerinspire.com/a11f-11001792/W1Z70FGYXPT140876/3957961676025655f7ba812adee58639.jpg, https://vehicle-images.dealerinspire.com/87bc-11001792/W1Z70FGYXPT140876/262ed58ca10dc3942679f3c2b0516c6e.jpg, https://vehicle-images.dealerinspire.com/35cf-11001792/W1Z70FGYXPT140876/1b5263dba46c8da94bc641afe8a08be8.jpg, https://vehicle-images.dealerinspire.com/8482-11001792/W1Z70FGYXPT140876/b28ed4c7b94eee136a8076374620657d.jpg, https://vehicle-images.dealerinspire.com/b53a-11001792/W1Z70FGYXPT140876/b4feb999be03c4f27a773323766e9742.jpg, https://vehicle-images.dealerinspire.com/26b4-11001792/W1Z70FGYXPT140876/06f430a4497918cc7b1f6d6a05edd9e0.jpg, https://vehicle-images.dealerinspire.com/c5ad-11001792/W1Z70FGYXPT140876/eb6954f80003c9897183eac540905e12.jpg, https://vehicle-images.dealerinspire.com/f6aa-11001792/W1Z70FGYXPT140876/9ac3da778df4967e3c982612b137d217.jpg, https://vehicle-images.dealerinspire.com/b205-11001792/W1Z70FGYXPT140876/66f7b899e9e312d43634be8569cb7a43.jpg, https://vehicle-images.dealerinspire.com/a09c-11001792/W1Z70FGYXPT140876/ec2f155b102be5178254d9711f85ddeb.jpg, https://vehicle-images.dealerinspire.com/5e20-11001792/W1Z70FGYXPT140876/97f1df47b99a8dd5657df9504703e763.jpg, https://vehicle-images.dealerinspire.com/4840-11001792/W1Z70FGYXPT140876/2adba1bfa4529dd3540cd4da633d8dd2.jpg'
$ws.Range("L5").Value = 'https://invimg2.autofunds.net/InventoryImages/2025/05/24/4181_2606707_65915339_4532424532025.jpg, https://invimg2.autofunds.net/InventoryImages/2025/05/24/4181_2606707_19641376_4523324522025.jpg, https://invimg2.autofunds.net/InventoryImages/2025/05/24/4181_2606707_62461927_4522624522025.jpg, https://invimg2.autofunds.net/InventoryImages/2025/05/24/4181_2606707_18467678_4522524522025.jpg, https://invimg2.autofunds.net/InventoryImages/2025/05/24/4181_2606707_21677126_4522624522025.jpg, https://invimg2.autofunds.net/InventoryImages/2025/05/24/4181_2606707_36118680_4472324472025.jpg, https://invimg2.autofunds.net/InventoryImages/2025/05/24/4181_2606707_97869350_4532424532025.jpg, https://invimg2.autofunds.net/InventoryImages/2025/05/24/4181_2606707_18009186_4472324472025.jpg, https://invimg2.autofunds.net/InventoryImages/2025/05/24/4181_2606707_548057_4532424532025.jpg, https://invimg2.autofunds.net/InventoryImages/2025/05/24/4181_2606707_33281711_4472324472025.jpg, https://invimg2.autofunds.net/InventoryImages/2025/05/24/4181_2606707_12604981_4472324472025.jpg, https://invimg2.autofunds.net/InventoryImages/2025/05/24/4181_2606707_24566635_4512424512025.jpg, https://invimg2.autofunds.net/InventoryImages/2025/05/24/4181_2606707_36912557_4512224512025.jpg, https://invimg2.autofunds.net/InventoryImages/2025/05/24/4181_2606707_12371612_4512724512025.jpg, https://invimg2.autofunds.net/InventoryImages/2025/05/24/4181_2606707_85840617_4522624522025.jpg, https://invimg2.autofunds.net/InventoryImages/2025/05/24/4181_2606707_72352111_4522624522025.jpg, https://invimg2.autofunds.net/InventoryImages/2025/05/24/4181_2606707_62834643_4522624522025.jpg, https://invimg2.autofunds.net/InventoryImages/2025/05/24/4181_2606707_17360961_4522624522025.jpg, https://invimg2.autofunds.net/InventoryImages/2025/05/24/4181_2606707_10595468_4522624522025.jpg, https://invimg2.autofunds.net/InventoryImages/2025/05/24/4181_2606707_60591147_4522624522025.jpg, https://invimg2.autofunds.net/InventoryImages/2025/05/24/4181_2606707_17256544_4482424482025.jpg, https://invimg2.autofunds.net/InventoryImages/2025/05/24/4181_2606707_10944751_4482724482025.jpg, https://invimg2.autofunds.net/InventoryImages/2025/05/24/4181_2606707_28592041_4522624522025.jpg, https://invimg2.autofunds.net/InventoryImages/2025/05/24/4181_2606707_18703302_4482724482025.jpg, https://invimg2.autofunds.net/InventoryImages/2025/05/24/4181_2606707_31226475_4482824482025.jpg, https://invimg2.autofunds.net/InventoryImages/2025/05/24/4181_2606707_18110232_4512524512025.jpg, https://invimg2.autofunds.net/InventoryImages/2025/05/24/4181_2606707_21110799_4482424482025.jpg, https://invimg2.autofunds.net/InventoryImages/2025/05/24/4181_2606707_61366732_4482924482025.jpg, https://invimg2.autofunds.net/InventoryImages/2025/05/24/4181_2606707_89139530_4472324472025.jpg, https://invimg2.autofunds.net/InventoryImages/2025/05/24/4181_2606707_88323496_4472324472025.jpg, https://invimg2.autofunds.net/InventoryImages/2025/05/24/4181_2606707_13694200_4461724462025.jpg, https://invimg2.autofunds.net/InventoryImages/2025/05/24/4181_2606707_52166234_4473024472025.jpg, https://invimg2.autofunds.net/InventoryImages/2025/05/24/4181_2606707_13629800_4472324472025.jpg, https://invimg2.autofunds.net/InventoryImages/2025/05/24/4181_2606707_16403809_4533124532025.jpg, https://invimg2.autofunds.net/InventoryImages/2025/05/24/4181_2606707_94431168_4472324472025.jpg, https://invimg2.autofunds.net/InventoryImages/2025/05/24/4181_2606707_10228949_4512224512025.jpg, https://invimg2.autofunds.net/InventoryImages/2025/05/24/4181_2606707_17577166_4532424532025.jpg, https://invimg2.autofunds.net/InventoryImages/2025/05/24/4181_2606707_20697111_4472324472025.jpg'

$ws.Range("L5").Style = "Normal"

$ws.Range("J9").Select()
